$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (testing person 1 - invitation 60107) - invitation opened
$ws.Range("H2").Value = "2015-04-19 17:07"

# Row 3 (testing person 2 - invitation 60107) - invitation opened
$ws.Range("H3").Value = "2015-04-19 17:07"

# Row 4 (the Moskovitzes - invitation 42652)
$ws.Range("G4").Value = 1
$ws.Range("H4").Value = "2015-04-19 21:41"
$ws.Range("I4").Value = ""

# Row 5 (army friend - invitation 20349)
$ws.Range("E5").Value = "Yes"
$ws.Range("G5").Value = 1
$ws.Range("H5").Value = "2015-04-19 19:01"
$ws.Range("I5").Value = "Vegan"

# Row 6 (Guest of army friend - invitation 20349)
$ws.Range("G6").Value = 1
$ws.Range("H6").Value = "2015-04-19 19:01"
